$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Assert-CellText($row, $expected) {
    $actual = $ws.Cells.Item($row, 1).Value2
    if ($actual -ne $expected) {
        throw "Row $row expected '$expected' but found '$actual'"
    }
}

# ---------------------------------------------------------------------------
# The sheet is a simple A-column report with three sections, each headed by
# a bold label and separated by a blank row:
#   1  Payroll Warnings                (title)
#   2  (blank)
#   3  Equip Log Missing Cost Code     (section header)
#   4..9  warning lines
#   10 (blank)
#   11 Equipment log entry with no matching time card entry (section header)
#   12..23 warning lines
#   24 (blank)
#   25 Equipment log with no operator  (section header)
#   26 warning line
#
# The edit: add a new equip-1009/job-225010/2026-02-13 pair of warnings to
# the first section, drop three now-stale 2026-02-13 duplicate entries from
# the second section, and add the matching "no operator" warning line for
# that same equip/job/date to the third section.
#
# Rows are processed from the bottom of the sheet upward so every row
# number used below always refers to the *current* (not-yet-shifted) sheet.
# ---------------------------------------------------------------------------

# 3) "Equipment log with no operator" section: insert the new equip-1009 line
#    right after the header, ahead of the existing equip-1042 line.
Assert-CellText 25 "Equipment log with no operator"
Assert-CellText 26 "- Equip #: 1042 Job: 225010, Date: 2026-02-10"

$ws.Rows(26).Insert()
$newCell = $ws.Cells.Item(26, 1)
$newCell.Value = "- Equip #: 1009 Job: 225010, Date: 2026-02-13"
$newCell.Style = "Normal"

# 2) "Equipment log entry with no matching time card entry" section: remove
#    the three stale 2026-02-13 duplicate lines (delete bottom-up so earlier
#    row numbers stay valid).
Assert-CellText 20 "- Luis Espinoza2250102026-02-13200/500"
$ws.Rows(20).Delete()

Assert-CellText 17 "- Luis Espinoza2250102026-02-13200/500"
$ws.Rows(17).Delete()

Assert-CellText 14 "- Gilberto Ortiz2250102026-02-13200/500"
$ws.Rows(14).Delete()

Assert-CellText 12 "- Agustin Avila2250102026-02-13200/500"
$ws.Rows(12).Delete()

# 1) "Equip Log Missing Cost Code" section: insert the two new equip-1009
#    warning lines right after the header.
Assert-CellText 3 "Equip Log Missing Cost Code"
Assert-CellText 4 "-  225010 2026-02-10 1042"

$ws.Rows(4).Insert()
$ws.Rows(4).Insert()

$c1 = $ws.Cells.Item(4, 1)
$c1.Value = "- Gilberto Ortiz 225010 2026-02-13 1009"
$c1.Style = "Normal"

$c2 = $ws.Cells.Item(5, 1)
$c2.Value = "-  225010 2026-02-13 1009"
$c2.Style = "Normal"
